$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Cells whose type flips between text placeholder ("0"/"***.*") and a real number ---
# Value is set first (quoted when text is required), then format is pasted from a
# same-row donor cell that already carries the destination style so the style index
# collapses back onto the pre-existing one instead of minting a new one.
$ws.Range("F14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D16").Value = 1
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = 100
$ws.Range("H16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$ws.Range("H26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Plain numeric value updates (no style/type change) ---
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 500
$ws.Range("M15").Value = 500
$ws.Range("N15").Value = 500
$ws.Range("C16").Value = 2
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 62.5
$ws.Range("L16").Value = -27.777777777777
$ws.Range("M16").Value = -45.833333333333
$ws.Range("N16").Value = -88.596491228070
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 27
$ws.Range("K17").Value = -14.814814814814
$ws.Range("L17").Value = -20.689655172413
$ws.Range("M17").Value = 76.923076923076
$ws.Range("N17").Value = -34.285714285714
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -46.666666666666
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = -10
$ws.Range("L18").Value = -14.285714285714
$ws.Range("M18").Value = -69.491525423728
$ws.Range("N18").Value = -93.155893536121
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -17.021276595744
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 82
$ws.Range("K19").Value = -10.975609756097
$ws.Range("L19").Value = 5.797101449275
$ws.Range("M19").Value = 32.727272727272
$ws.Range("N19").Value = -20.652173913043
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 20
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = 52.631578947368
$ws.Range("L20").Value = 141.666666666667
$ws.Range("M20").Value = 45
$ws.Range("N20").Value = -88.537549407114
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -4.761904761904
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = -7.216494845360
$ws.Range("I21").Value = 164
$ws.Range("J21").Value = 159
$ws.Range("K21").Value = 3.144654088050
$ws.Range("L21").Value = 9.333333333333
$ws.Range("M21").Value = -4.651162790697
$ws.Range("N21").Value = -78.421052631578
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -22.580645161290
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 13.095238095238
$ws.Range("I24").Value = 136
$ws.Range("J24").Value = 151
$ws.Range("K24").Value = -9.933774834437
$ws.Range("L24").Value = 1.492537313432
$ws.Range("M24").Value = -1.449275362318
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 55
$ws.Range("I25").Value = 50
$ws.Range("J25").Value = 35
$ws.Range("K25").Value = 42.857142857142
$ws.Range("L25").Value = 28.205128205128
$ws.Range("M25").Value = -1.960784313725
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = 133.333333333333
$ws.Range("L26").Value = 133.333333333333
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 75
$ws.Range("I27").Value = 15
$ws.Range("J27").Value = 5
$ws.Range("L27").Value = 66.666666666666
$ws.Range("L30").Value = 0
